$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (pushes old rows 11-13 down to 12-14).
# Copy row 10's formatting (style pattern A2,B2,C2,D3,E2,F2,G2,H2,I5) which
# matches the formatting the new row needs, so the per-cell style indices
# carry over correctly.
$ws.Rows.Item(10).Copy()
$ws.Rows.Item(11).Insert(-4121)

# Populate the new row 11 content.
$ws.Range("A11").Value = "albert"
$ws.Range("E11").Value = "Nebraska"
$ws.Range("G11").ClearContents()
$ws.Range("I11").Value = "philbert in Nebraska as described by Anon."

# Set the row height for the new row.
$ws.Rows.Item(11).RowHeight = 75

# Update the window view / selection to match the author's edit session.
$ws.Range("I11").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 35660
$excel.ActiveWindow.Top = 6160
